$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text, matching the source data
# (these are formatted numeric-looking strings, e.g. "599.71" or "65.691.84",
# which Excel would otherwise auto-convert to numbers on assignment).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '65.691.84'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.675.30'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '599.71'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').Value = '156.61'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  +5.37%  '
$ws.Range('D9').Value = '0.130'
$ws.Range('E9').Value = '  +4.96%  '
$ws.Range('D10').Value = '0.400'
$ws.Range('E10').Value = '  -0.44%  '
$ws.Range('D11').Value = '5.87'
$ws.Range('E11').Value = '  -2.65%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = '29.31'
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').Value = '0.0000197'
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').Value = '3.155.40'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '66.334.81'
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').Value = '2.678.60'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '12.90'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('D19').Value = '4.79'
$ws.Range('E19').Value = '  -1.94%  '
$ws.Range('D20').Value = '7.55'
$ws.Range('E20').Value = '  +0.23%  '
$ws.Range('D21').Value = '351.97'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '69.82'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('E24').Value = '  +4.34%  '
$ws.Range('D25').Value = '9.62'
$ws.Range('E25').Value = '  -2.50%  '
$ws.Range('D26').Value = '1.64'
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').Value = '0.166'
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('E28').Value = '  -5.38%  '
$ws.Range('D29').Value = '8.03'
$ws.Range('E29').Value = '  -3.56%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.13'
$ws.Range('E31').Value = '  -3.27%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '526.89'
$ws.Range('E32').Value = '  -3.22%  '
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D34').Value = '6.45'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').Value = '5.46'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').Value = '0.423'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '20.62'
$ws.Range('E37').Value = '  -0.92%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').Value = '157.94'
$ws.Range('E39').Value = '  -3.21%  '
$ws.Range('D40').Value = '1.93'
$ws.Range('E40').Value = '  -2.79%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = '163.98'
$ws.Range('E42').Value = '  -4.23%  '
$ws.Range('D43').Value = '4.12'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('D44').Value = '2.30'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').Value = '0.0609'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = '22.79'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').Value = '0.640'
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('D49').Value = '0.0₆0264'
$ws.Range('E49').Value = '  +15.94%  '
$ws.Range('E50').Value = '  +0.91%  '
$ws.Range('D51').Value = '20.13'
$ws.Range('E51').Value = '  -4.21%  '

# Drop the temporary text-number-format so the cells keep the original
# (default) style, same as before the edit.
$ws.Range("D2:D51").ClearFormats()

Write-Host "Update complete"
